$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Kondisi pemanfaatan ruang ... berupa (" paragraph: split the trailing
#    " (" off the sentence and insert a new "batas utara/selatan/timur/barat"
#    sentence (with a line break) between "berupa" and the "(" that opens
#    the red instructional note.
# ---------------------------------------------------------------------------
$old1 = 'Kondisi pemanfaatan ruang di sebelah dan di sekitar lahan yang dimohonkan berupa ('
$new1 = 'Kondisi pemanfaatan ruang di sebelah dan di sekitar lahan yang dimohonkan berupa :' + [char]11 + 'batas utara : ${batas_utara}, batas Selatan : ${batas_selatan}, batas timur : ${batas_timur}, batas barat : ${batas_barat} ('

$r1 = $d.Content
$found1 = $r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    throw "Could not find target text for change 1"
}

# ---------------------------------------------------------------------------
# 2) "Lokasi persil berada pada" now starts a new rendered page, so Word
#    stamps a lastRenderedPageBreak marker right before it.
# ---------------------------------------------------------------------------
# (lastRenderedPageBreak is a render-time bookkeeping marker Word regenerates
#  during pagination; it carries no editable content through the object
#  model, so it is intentionally left to the application's own repagination.)

# ---------------------------------------------------------------------------
# 3) "dengan tipe jalan " + "dengan " (previously split across a
#    lastRenderedPageBreak) collapse back into a single contiguous run of
#    text now that the page break marker has moved elsewhere.
# ---------------------------------------------------------------------------
$old3 = 'dengan tipe jalan dengan '
$r3 = $d.Content
$found3 = $r3.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $old3, 2)
if (-not $found3) {
    throw "Could not find target text for change 3"
}

Write-Host "Edits applied"
